$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Production (MW)" values for rows 2..30 (a new day's data was added to the
# rolling window). Rows 31..97 keep their existing value of 0 (only the timestamp shifts).
$newValuesB = @(2354, 2351, 2320, 2284, 2277, 2253, 2246, 2250, 2234, 2236, 2202, 2141, 2095, 2063, 2013, 1967, 1934, 1892, 1795, 1745, 1761, 1734, 1679, 1680, 1706, 1731, 1754, 1758, 0)

for ($r = 2; $r -le 97; $r++) {
    # Shift every timestamp in column A forward by one day.
    $cellA = $ws.Cells.Item($r, 1)
    $oldDate = $cellA.Value()
    $cellA.Value = $oldDate.AddDays(1)

    # Update column B where a new value is supplied (rows 2..30); otherwise leave as-is.
    $idx = $r - 2
    if ($idx -lt $newValuesB.Length) {
        $ws.Cells.Item($r, 2).Value = $newValuesB[$idx]
    }
}
